# Update "想去人数" (interested-count) figures and mark the cancelled
# 乐平·CY境界次元第三届动漫游戏庆典 event on both the "展览" sheet and the
# aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    if ($sheetName -eq "展览") {
        $ws.Range("F2").Value = 1489
        $ws.Range("F4").Value = 1790
        $ws.Range("F10").Value = 556
        $ws.Range("F13").Value = 83

        $ws.Range("C17").Value = "乐平·CY境界次元第三届动漫游戏庆典（取消）"
        $ws.Range("G17").Value = "不可售"

        $ws.Range("F19").Value = 5005
        $ws.Range("F20").Value = 52
        $ws.Range("F21").Value = 833
        $ws.Range("F22").Value = 114
        $ws.Range("F23").Value = 2251
        $ws.Range("F26").Value = 2096
    }
    else {
        $ws.Range("F2").Value = 1489
        $ws.Range("F4").Value = 1790
        $ws.Range("F10").Value = 556
        $ws.Range("F13").Value = 83

        $ws.Range("C17").Value = "乐平·CY境界次元第三届动漫游戏庆典（取消）"
        $ws.Range("G17").Value = "不可售"

        $ws.Range("F19").Value = 5005
        $ws.Range("F21").Value = 52
        $ws.Range("F23").Value = 833
        $ws.Range("F24").Value = 114
        $ws.Range("F25").Value = 2251
        $ws.Range("F28").Value = 2096
    }
}
